$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New equal weights (1/7) for the "Initial Weights" column (B2:B8)
$ws.Range("B2").Value = 0.1428571428571428
$ws.Range("B3").Value = 0.1428571428571428
$ws.Range("B4").Value = 0.1428571428571428
$ws.Range("B5").Value = 0.1428571428571428
$ws.Range("B6").Value = 0.1428571428571428
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("B8").Value = 0.1428571428571428

# Updated "Opt Portfolio" column (C2:C8)
$ws.Range("C2").Value = 0.1535818480114002
$ws.Range("C3").Value = 0.1522119943432627
$ws.Range("C4").Value = 0.1315750976089232
$ws.Range("C5").Value = 0.1316127775234491
$ws.Range("C6").Value = 0.1315393276269073
$ws.Range("C7").Value = 0.1495287225222207
$ws.Range("C8").Value = 0.1499502323638368

# Updated "Opt Portfolio with View" column (D2:D8)
$ws.Range("D2").Value = 0.1535818492496737
$ws.Range("D3").Value = 0.1522119928421055
$ws.Range("D4").Value = 0.1315750984268254
$ws.Range("D5").Value = 0.1316127770923149
$ws.Range("D6").Value = 0.1315393269293997
$ws.Range("D7").Value = 0.1495287223521467
$ws.Range("D8").Value = 0.1499502331075342
